$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2 and 3 with the new control-point values
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 12
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 3

# Remove the now-obsolete last rows (4 and 5)
$ws.Range("A4:B5").Delete()
